$d = $word.ActiveDocument

$pairs = @(
    @("80×50=4000", "33×99=3267"),
    @("88×63=5544", "88×67=5896"),
    @("73×49=3577", "75×62=4650"),
    @("60×78=4680", "70×44=3080"),
    @("30×92=2760", "64×96=6144"),
    @("45×78=3510", "75×53=3975"),
    @("25×76=1900", "81×80=6480"),
    @("13×42=546",  "61×23=1403"),
    @("30×60=1800", "38×75=2850"),
    @("17×27=459",  "13×82=1066"),
    @("58×85=4930", "55×85=4675"),
    @("47×27=1269", "57×15=855"),
    @("87×59=5133", "35×49=1715"),
    @("97×58=5626", "67×91=6097"),
    @("16×73=1168", "43×81=3483"),
    @("62×88=5456", "46×17=782"),
    @("30×19=570",  "66×60=3960"),
    @("52×91=4732", "25×57=1425"),
    @("35×46=1610", "99×26=2574"),
    @("95×87=8265", "38×84=3192"),
    @("40×70=2800", "19×32=608"),
    @("29×26=754",  "78×17=1326"),
    @("27×79=2133", "31×34=1054"),
    @("86×50=4300", "30×79=2370"),
    @("91×19=1729", "85×56=4760")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
